$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (row 1): "from" sources, columns B..K
$fromHeaders = @(
    "Q_from_net1",
    "Q_from_net2",
    "Q_from_CHP1",
    "Q_from_CHP2",
    "Q_from_solar_th1",
    "Q_from_solar_th2",
    "Q_from_pvt1",
    "Q_from_pvt2",
    "Q_from_heat_pump1",
    "Q_from_heat_pump2"
)

for ($i = 0; $i -lt $fromHeaders.Length; $i++) {
    $col = 2 + $i  # B = 2
    $ws.Cells.Item(1, $col).Value = $fromHeaders[$i]
}

# Row headers (column A): "to" destinations, rows 2..7
$toHeaders = @(
    "param_Q_to_demand1",
    "param_Q_to_demand2",
    "Q_to_net1",
    "Q_to_net2",
    "Q_to_heat_pump1",
    "Q_to_heat_pump2"
)

for ($i = 0; $i -lt $toHeaders.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $toHeaders[$i]
}

# Matrix body values B2:K7 -- "to_from" style labels for rows 2-5, numeric 0 for rows 6-7
# and numeric 0 in B4:C5 (net-to-net entries)
$toKeys = @("demand1", "demand2", "net1", "net2")
$fromKeys = @("net1", "net2", "CHP1", "CHP2", "solar_th1", "solar_th2", "pvt1", "pvt2", "heat_pump1", "heat_pump2")

$netKeys = @("net1", "net2")

for ($r = 0; $r -lt 4; $r++) {
    $rowNum = 2 + $r
    for ($c = 0; $c -lt $fromKeys.Length; $c++) {
        $colNum = 2 + $c
        $toName = $toKeys[$r]
        $fromName = $fromKeys[$c]

        if (($netKeys -contains $fromName) -and ($netKeys -contains $toName)) {
            # net-to-net entries (B4:C5) stay numeric 0 (no net-to-net flow string)
            $ws.Cells.Item($rowNum, $colNum).Value = 0
        } else {
            $ws.Cells.Item($rowNum, $colNum).Value = "Q_" + $fromName + "_" + $toName
        }
    }
}

# Rows 6 and 7 (Q_to_heat_pump1 / Q_to_heat_pump2) are all zeros across B..K
for ($row = 6; $row -le 7; $row++) {
    for ($col = 2; $col -le 11; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}

# Apply bold/border/center formatting to the newly added header cells (C1:K1, A4:A7)
# by copying the format already present on B1 (which matches the style used by
# A2/A3/B1 in the original workbook), avoiding creation of extra/unused style entries.
# (Each contiguous area is pasted separately since multi-area destination ranges
# only apply PasteSpecial to the first area.)
$src = $ws.Range("B1")
$src.Copy()
$ws.Range("C1:K1").PasteSpecial(-4122)
$src.Copy()
$ws.Range("A4:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
